# Apply updates to the "year" sheet of the EMEP/NFR09 scaling mapping workbook:
# add a new scaling-year rule row for "lux" (Luxembourg) / RoadRail sector.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("year")
$ws.Activate()

$ws.Cells.Item(4, 1).Value = "lux"
$ws.Cells.Item(4, 2).Value = "RoadRail"
$ws.Cells.Item(4, 3).Value = "NA"
$ws.Cells.Item(4, 4).Value = "NA"
$ws.Cells.Item(4, 5).Value = "NA"
$ws.Cells.Item(4, 6).Value = 1990
$ws.Cells.Item(4, 7).Value = 2020
$ws.Cells.Item(4, 8).Value = "NA"
$ws.Cells.Item(4, 9).Value = "Avoid imlied Nox EF dip 1986-1989"

# Match the new active selection on the freshly-added row.
$ws.Range("A4:XFD4").Select()

# Reposition the workbook window, as recorded by Excel when the file was last saved.
$win = $wb.Windows.Item(1)
$win.Left = 6680
$win.Top = 2280
